$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data rows 2-18 (A:H) with the new dataset values
$data = @(
    @(0, "falling", -0.3275403976440429, -1.032773733139038, -0.2356588244438171, 1.223868489265442, 1.015716195106506, 0.2379320114850998),
    @(100, "falling", -1.590312331914902, -0.9565697163343428, -1.569369990378618, -0.2808452844619751, -0.5580254197120667, 0.9086620807647704),
    @(200, "falling", 0.2956193089485168, -0.9267413020133971, -1.241599485278129, -0.4292855560779571, -1.047786593437195, 0.4741841554641723),
    @(300, "falling", -0.6574213504791258, -0.9026327282190323, -1.075559064745903, 0.2979495227336883, -0.4011857509613037, 0.1365283876657486),
    @(400, "falling", 0.3167376518249511, -0.9883218407630923, -1.371818482875824, -0.0152716310694813, -0.2217440903186798, -0.104763388633728),
    @(500, "falling", 0.4182748794555663, -0.9399109184741974, -1.527341216802597, -0.0087048299610614, 0.0255036242306232, -0.0922406539320945),
    @(600, "falling", 0.2231501340866089, -0.9057424068450929, -1.459591150283814, 0.0565050356090068, -0.0242818929255008, 0.0531452745199203),
    @(700, "falling", -0.4751685261726377, -0.88183431327343, -1.828994989395142, 0.0300851128995418, 0.0068722339347004, -0.0010690141934901),
    @(800, "falling", -0.972740650177002, -1.058116793632508, -2.322797894477844, 0.1162171140313148, -0.3998112976551056, 0.1134682223200798),
    @(900, "falling", -1.244342982769013, -1.110738858580589, -2.960273459553719, -0.1954768747091293, -0.8119926452636719, 0.1353066563606262),
    @(1000, "falling", -1.067278563976287, -1.54173904657364, -3.944739699363709, -0.2617557644844055, -0.7904596328735352, -0.3645338416099548),
    @(1100, "falling", -0.4881105124950414, -2.259673044085503, -4.275467619299889, 0.3368921875953674, 0.3489567637443542, 1.531897306442261),
    @(1200, "falling", 0.24235272407532, -1.947445154190064, -1.148832976818081, -0.4462370574474334, 2.508059978485107, 0.2651155292987823),
    @(1300, "falling", 6.217258334159848, -2.320076808333399, 9.556089490652074, 0.6637051105499268, -0.5948300361633301, -0.4066835343837738),
    @(1400, "falling", 10.25036442279817, 0.2667215764522569, 19.65825939178468, 0.1577559560537338, -5.182427883148193, -2.026240110397339),
    @(1500, "falling", -5.869112968444829, -2.593451023101809, -7.71562802791596, -3.560580730438232, -1.920713067054748, 1.063363671302795),
    @(1600, "falling", 14.12184143066406, 8.388177871704102, 14.74949073791504, -6.374378681182861, -4.941136360168457, -10.62203025817871)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
}
